$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 (header-ish values) tweaks
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2 (CON) tweaks
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = -0.60000000000000009
$ws.Range("D2").Value = -0.8
$ws.Range("E2").Value = -0.05

# Row 3 (STR) tweaks
$ws.Range("B3").Value = -1.2000000000000002
$ws.Range("C3").Value = -0.1
$ws.Range("D3").Value = -2.0500000000000003
$ws.Range("E3").Value = -2.15

# Update the selection to the edited block, matching the author's selection change
$excel.Goto($ws.Range("B1:E3"))
